$d = $word.ActiveDocument

# The footer block at the end of the document currently looks like:
#   ... "LOQ4031: Química Geral I (Requisito fraco)"
#   (empty paragraph)
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
#   (empty paragraph)
#   (page-break paragraph)
#
# The edit removes the empty paragraph plus the two "Ver no Jupiter..." /
# "© 2020 ..." paragraphs entirely, so the requisito line is immediately
# followed by the trailing empty paragraph + page-break paragraph.

# Locate the paragraph that keeps the "LOQ4031" requisito line (content to
# keep) and the paragraph holding the "Creative Commons" copyright notice
# (last paragraph to remove).
$keepIdx = -1
$lastRemoveIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($keepIdx -eq -1 -and $text -like "*LOQ4031*") {
        $keepIdx = $i
    }
    if ($text -like "*Creative Commons*") {
        $lastRemoveIdx = $i
    }
}

if ($keepIdx -ne -1 -and $lastRemoveIdx -ne -1 -and $lastRemoveIdx -gt $keepIdx) {
    $startDel = $d.Paragraphs.Item($keepIdx + 1).Range.Start
    $endDel = $d.Paragraphs.Item($lastRemoveIdx).Range.End

    $delRange = $d.Range($startDel, $endDel)
    $delRange.Delete()
}
